$wb = $excel.ActiveWorkbook

$msg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# --- Sheet "Bico" (column H = Obs -> Obs_relatorio, new column I = Obs_sped) ---
$bico = $wb.Worksheets.Item("Bico")

$bico.Range("H1").Value = "Obs_relatorio"
$bico.Range("I1").Value = "Obs_sped"

for ($r = 2; $r -le 13; $r++) {
    $bico.Cells.Item($r, 8).Value = $msg
    $bico.Cells.Item($r, 9).Value = ""
}

# --- Sheet "Tanque" (column F = Obs -> Obs_relatorio, new column G = Obs_sped) ---
$tanque = $wb.Worksheets.Item("Tanque")

$tanque.Range("F1").Value = "Obs_relatorio"
$tanque.Range("G1").Value = "Obs_sped"

for ($r = 2; $r -le 6; $r++) {
    $tanque.Cells.Item($r, 6).Value = $msg
    $tanque.Cells.Item($r, 7).Value = ""
}
